# Generate Report for Handback
# Replaces the file identifiers / timestamps for the two handback rows across
# all three worksheets (Overview, zh-cn, de-de), and updates the matching
# hyperlink display text.

$wb = $excel.ActiveWorkbook

$oldGuid1 = "32b7cbc2-45b2-451f-8e8d-1ab2cc0211a5"
$newGuid1 = "d5e47a05-7745-416b-a80f-2a61e8b68b7a"
$oldGuid2 = "e763c5cb-5d0e-42a3-a729-ac235ed9e9fd"
$newGuid2 = "ffffa777b6d6-87e7-4a44-a1ec-7a147d912d17"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = "$newGuid1.md"
$ws.Range("B2").Value = "e2e\$newGuid1.md"
$ws.Range("G2").Value = "2016-08-18 23:04:20"

$ws.Range("A3").Value = "$newGuid2.md"
$ws.Range("B3").Value = "e2e\$newGuid2.md"
$ws.Range("G3").Value = "2016-08-18 23:04:20"

# The COM shim only ever appends new hyperlink entries rather than updating
# them in place, so clear the sheet's hyperlinks and recreate them (this
# reuses the same r:id order/values as before). The underlying hyperlink
# target URLs are not changed by this edit - only their displayed text is.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9ea7d55a6eb73e723419466ffb2cc5f4b451d623/e2e/$oldGuid1.md", "", "", "e2e\$newGuid1.md")
$ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9ea7d55a6eb73e723419466ffb2cc5f4b451d623/e2e/$oldGuid2.md", "", "", "e2e\$newGuid2.md")

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2").Value = "$newGuid1.md"
$ws.Range("G2").Value = "$newGuid1.0fb762ba73f259f3995538a124d3941c824ea16e.zh-cn.xlf"
$ws.Range("H2").Value = "2016-08-18 23:04:14"
$ws.Range("I2").Value = "$newGuid1.md"
$ws.Range("J2").Value = "$newGuid1.0fb762ba73f259f3995538a124d3941c824ea16e.zh-cn.xlf"
$ws.Range("K2").Value = "2016-08-18 23:04:32"

$ws.Range("A3").Value = "$newGuid2.md"
$ws.Range("G3").Value = "$newGuid1.0fb762ba73f259f3995538a124d3941c824ea16e.zh-cn.xlf"
$ws.Range("H3").Value = "2016-08-18 23:04:14"
$ws.Range("I3").Value = "$newGuid2.md"
$ws.Range("J3").Value = "$newGuid1.0fb762ba73f259f3995538a124d3941c824ea16e.zh-cn.xlf"
$ws.Range("K3").Value = "2016-08-18 23:04:32"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9ea7d55a6eb73e723419466ffb2cc5f4b451d623/e2e/$oldGuid1.md", "", "", "$newGuid1.md")
$ws.Hyperlinks.Add($ws.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/a1dbadd792794bfdea19d08e694fb56e20e712f5/e2e/$oldGuid1.md", "", "", "$newGuid1.md")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9ea7d55a6eb73e723419466ffb2cc5f4b451d623/e2e/$oldGuid2.md", "", "", "$newGuid2.md")
$ws.Hyperlinks.Add($ws.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/a1dbadd792794bfdea19d08e694fb56e20e712f5/e2e/$oldGuid2.md", "", "", "$newGuid2.md")

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2").Value = "$newGuid1.md"
$ws.Range("G2").Value = "$newGuid1.0fb762ba73f259f3995538a124d3941c824ea16e.de-de.xlf"
$ws.Range("H2").Value = "2016-08-18 23:04:20"
$ws.Range("I2").Value = "$newGuid1.md"
$ws.Range("J2").Value = "$newGuid1.0fb762ba73f259f3995538a124d3941c824ea16e.de-de.xlf"
$ws.Range("K2").Value = "2016-08-18 23:04:39"

$ws.Range("A3").Value = "$newGuid2.md"
$ws.Range("G3").Value = "$newGuid1.0fb762ba73f259f3995538a124d3941c824ea16e.de-de.xlf"
$ws.Range("H3").Value = "2016-08-18 23:04:20"
$ws.Range("I3").Value = "$newGuid2.md"
$ws.Range("J3").Value = "$newGuid1.0fb762ba73f259f3995538a124d3941c824ea16e.de-de.xlf"
$ws.Range("K3").Value = "2016-08-18 23:04:39"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9ea7d55a6eb73e723419466ffb2cc5f4b451d623/e2e/$oldGuid1.md", "", "", "$newGuid1.md")
$ws.Hyperlinks.Add($ws.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/2fd83b63ba3485a4945929b8d546b47074c2b8b2/e2e/$oldGuid1.md", "", "", "$newGuid1.md")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9ea7d55a6eb73e723419466ffb2cc5f4b451d623/e2e/$oldGuid2.md", "", "", "$newGuid2.md")
$ws.Hyperlinks.Add($ws.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/2fd83b63ba3485a4945929b8d546b47074c2b8b2/e2e/$oldGuid2.md", "", "", "$newGuid2.md")

Write-Host "Handback status report regenerated."
